$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dosen")
$c = $ws.Range("C72")
$c.Font.Name = "Calibri Light"
Write-Host "interim s:" $c.Style.Name
